$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.977.44"
$ws.Range("E2").Value = "  +0.22%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.639.08"
$ws.Range("E3").Value = "  -0.43%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.06"
$ws.Range("E5").Value = "  -0.21%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5138"
$ws.Range("E6").Value = "  +0.61%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.002"
$ws.Range("E7").Value = "  -0.41%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2574"
$ws.Range("E8").Value = "  -0.11%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06354"
$ws.Range("E9").Value = "  -1.05%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.78"
$ws.Range("E10").Value = "  +0.31%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07769"
$ws.Range("E11").Value = "  -0.07%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.275"
$ws.Range("E12").Value = "  -0.84%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.632.89"
$ws.Range("E13").Value = "  -0.85%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5457"
$ws.Range("E14").Value = "  -0.37%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0₅7749"
$ws.Range("E15").Value = "  -1.92%  "
$ws.Range("E16").Value = "  -1.02%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "25.991.37"
$ws.Range("E17").Value = "  -0.05%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9976"
$ws.Range("E18").Value = "  -0.84%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "197.43"
$ws.Range("E19").Value = "  -0.20%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.440"
$ws.Range("E20").Value = "  +0.16%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.928"
$ws.Range("E21").Value = "  -1.14%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.084"
$ws.Range("E22").Value = "  +0.16%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.004"
$ws.Range("E23").Value = "  -0.41%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.935"
$ws.Range("E24").Value = "  +4.01%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "142.11"
$ws.Range("E25").Value = "  +1.01%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1235"
$ws.Range("E26").Value = "  +7.60%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.844"
$ws.Range("E27").Value = "  -0.79%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.63"
$ws.Range("E28").Value = "  -0.95%  "
$ws.Range("E29").Value = "  +0.04%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.04847"
$ws.Range("E30").Value = "  -3.39%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.284"
$ws.Range("E31").Value = "  +0.22%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.216"
$ws.Range("E32").Value = "  +0.43%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.538"
$ws.Range("E33").Value = "  -0.41%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.375"
$ws.Range("E34").Value = "  +0.34%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9138"
$ws.Range("E35").Value = "  +2.08%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.570"
$ws.Range("E36").Value = "  -0.75%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.5550"
$ws.Range("E37").Value = "  -0.15%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.103.65"
$ws.Range("E38").Value = "  -2.53%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01568"
$ws.Range("E39").Value = "  +0.09%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.001"
$ws.Range("E40").Value = "  -0.53%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.521"
$ws.Range("E41").Value = "  -1.78%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.559"
$ws.Range("E42").Value = "  -1.95%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8068"
$ws.Range("E43").Value = "  -1.04%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "99.28"
$ws.Range("E44").Value = "  -0.43%  "
$ws.Range("E45").Value = "  -2.82%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.781.05"
$ws.Range("E46").Value = "  -0.21%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4538"
$ws.Range("E47").Value = "  +0.02%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "55.08"
$ws.Range("E48").Value = "  -0.48%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.9993"
$ws.Range("E49").Value = "  -0.79%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05212"
$ws.Range("E50").Value = "  +2.38%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.487"
$ws.Range("E51").Value = "  +1.08%  "
